# Updated queries in the 6 failing UBC01 scripts:
# the StatQuery (column C) text on the CasesTab/SamplesTab/FilesTab rows
# is replaced with a new Cypher query that adds cohort/aliquot counting
# and parameterised filter clauses.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newStatQuery = @'
MATCH (p:program)<--(s:study)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
      WHERE (size([]) = 0 OR s.clinical_study_designation IN [])
        AND (s.study_disposition = 'Unrestricted')
        AND (size([]) = 0 OR s.clinical_study_type IN [])
        AND s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in ['T2N0M0', 'Not Applicable']
        AND (size([]) = 0 OR demo.sex IN [])
        AND (size([]) = 0 OR demo.neutered_indicator IN [])
        AND (size([]) = 0 OR diag.disease_term IN [])
        AND (size([]) = 0 OR diag.primary_disease_site IN [])
        AND (size([]) = 0 OR diag.stage_of_disease IN [])
        AND (size([]) = 0 OR diag.best_response IN [])
    OPTIONAL MATCH (c)-->(co:cohort)
    OPTIONAL MATCH (f:file)-[*]->(c)
    OPTIONAL MATCH (f)-->(parent)
    OPTIONAL MATCH (samp:sample)-->(c)
    OPTIONAL MATCH (samp)<--(al:aliquot)
    WITH DISTINCT c AS c, p, s, co, demo, diag, f, parent, samp, al
      WHERE (size([]) = 0 OR samp.summarized_sample_type IN [])
        AND (size([]) = 0 OR samp.specific_sample_pathology IN [])
        AND (size([]) = 0 OR samp.sample_site IN [])
        AND (size([]) = 0 OR head(labels(parent)) IN [])
        AND (size([]) = 0 OR f.file_type IN [])
        AND (size([]) = 0 OR f.file_format IN [])
    WITH c.case_id AS case_id,
         s.clinical_study_designation AS study_code,
         s.clinical_study_type AS study_type,
         co.cohort_description AS cohort,
         demo.breed AS breed,
         diag.disease_term AS diagnosis,
         diag.stage_of_disease AS stage_of_disease,
         diag.primary_disease_site AS disease_site,
         demo.patient_age_at_enrollment AS age,
         demo.sex AS sex,
         demo.neutered_indicator AS neutered_status,
         demo.weight AS weight,
         diag.best_response AS response_to_treatment,
         samp.sample_id AS sample_id,
         f.uuid AS file_id,
         al
    RETURN
COUNT(DISTINCT file_id) as number_of_files,
COUNT(DISTINCT sample_id) as number_of_sample,
COUNT(DISTINCT case_id) as number_of_cases,
COUNT(DISTINCT study_code) as number_of_study,
COUNT(DISTINCT al) as number_of_aliquot
'@
$newStatQuery = $newStatQuery.TrimEnd("`r", "`n")

# Replace the StatQuery column (C) for the 3 tab rows with the new query.
$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Row heights grow to the max (409.5pt, Excel's 409.5pt cap) to fit the
# now-longer wrapped text.
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# The view scrolled down (top-left visible cell became B4) and the
# selection moved to C4.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 2
$win.TopLeftCell = $ws.Range("B4")
$ws.Range("C4").Select()
